# add gender to doctors in import
#
# The "doctors" import sheet gains a new "gender"/"male" column, inserted
# right after "nameEn" (i.e. as the new column C). Every existing column
# from the old "specialty" column onward shifts one place to the right,
# and the last existing column ("homeTel" / its phone-number value) is
# dropped off the end since the sheet stays 8 columns wide (A:H).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- snapshot the current header row (1) and data row (2), columns C..G,
#     before anything gets overwritten. (.Value2 is used for reads; the
#     plain .Value getter does not round-trip scalars reliably here.)
$c1 = $ws.Range("C1").Value2
$d1 = $ws.Range("D1").Value2
$e1 = $ws.Range("E1").Value2
$f1 = $ws.Range("F1").Value2
$g1 = $ws.Range("G1").Value2

$c2 = $ws.Range("C2").Value2
$d2 = $ws.Range("D2").Value2
$e2 = $ws.Range("E2").Value2
$f2 = $ws.Range("F2").Value2
$g2 = $ws.Range("G2").Value2

# --- shift columns C..G right into D..H (H's old content -- "homeTel" /
#     404474444 -- is intentionally dropped, it falls off the end)
$ws.Range("H1").Value2 = $g1
$ws.Range("G1").Value2 = $f1
$ws.Range("F1").Value2 = $e1
$ws.Range("E1").Value2 = $d1
$ws.Range("D1").Value2 = $c1

$ws.Range("H2").Value2 = $g2
$ws.Range("G2").Value2 = $f2
$ws.Range("F2").Value2 = $e2
$ws.Range("E2").Value2 = $d2
$ws.Range("D2").Value2 = $c2

# --- populate the new gender column at C
$ws.Range("C1").Value2 = "gender"
$ws.Range("C2").Value2 = "male"

# match the author's final selection
$ws.Range("C3").Select()
